$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 299.22223
$ws.Range("I2").Value = 236.625
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 236.625
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -123.625
$ws.Range("N2").Value = -1026
$ws.Range("H12").Value = 283.92307
$ws.Range("I12").Value = 283.83334
$ws.Range("K12").Value = 283.83334
$ws.Range("M12").Value = -113.83334
$ws.Range("H17").Value = 2031.6
$ws.Range("J17").Value = 2031.6
$ws.Range("L17").Value = 6094.799999999999
$ws.Range("N17").Value = -6430.799999999999
$ws.Range("H18").Value = 1829.8462
$ws.Range("I18").Value = 1149
$ws.Range("K18").Value = 1149
$ws.Range("M18").Value = -865
$ws.Range("H28").Value = 42396
$ws.Range("I28").Value = 44152.78
$ws.Range("J28").Value = 1990
$ws.Range("K28").Value = 44152.78
$ws.Range("L28").Value = 1990
$ws.Range("M28").Value = -43667.78
$ws.Range("N28").Value = -2960
$ws.Range("H29").Value = 385.57144
$ws.Range("J29").Value = 766.6667
$ws.Range("L29").Value = 2300.0001
$ws.Range("N29").Value = -2862.0001
$ws.Range("H38").Value = 52.4375
$ws.Range("I38").Value = 52.4375
$ws.Range("K38").Value = 157.3125
$ws.Range("M38").Value = 214.6875
$ws.Range("H98").Value = 4399.8096
$ws.Range("I98").Value = 4182.5293
$ws.Range("K98").Value = 4182.5293
$ws.Range("M98").Value = -2684.5293
$ws.Range("H107").Value = 399.5
$ws.Range("I107").Value = 400.85715
$ws.Range("K107").Value = 400.85715
$ws.Range("M107").Value = 1519.14285
$ws.Range("H122").Value = 4399.8096
$ws.Range("I122").Value = 4182.5293
$ws.Range("K122").Value = 12547.5879
$ws.Range("M122").Value = -10097.5879
$ws.Range("H127").Value = 2033.3334
$ws.Range("I127").Value = 1000
$ws.Range("J127").Value = 2550
$ws.Range("K127").Value = 3000
$ws.Range("L127").Value = 7650
$ws.Range("M127").Value = 1960
$ws.Range("N127").Value = -17570
$ws.Range("H137").Value = 3303.7026
$ws.Range("I137").Value = 2302.48
$ws.Range("J137").Value = 5389.5835
$ws.Range("K137").Value = 6907.440000000001
$ws.Range("L137").Value = 16168.7505
$ws.Range("M137").Value = -4357.440000000001
$ws.Range("N137").Value = -21268.7505
$ws.Range("H141").Value = 8064.8
$ws.Range("I141").Value = 8081.1665
$ws.Range("J141").Value = 7999.3335
$ws.Range("K141").Value = 24243.4995
$ws.Range("L141").Value = 23998.0005
$ws.Range("M141").Value = -19063.4995
$ws.Range("N141").Value = -34358.00049999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6477.2246
$ws.Range("I32").Value = 5130.8374
$ws.Range("J32").Value = 18445.111
$ws.Range("K32").Value = 5130.8374
$ws.Range("L32").Value = 18445.111
$ws.Range("M32").Value = -4843.8374
$ws.Range("N32").Value = -19019.111
$ws.Range("H61").Value = 1597.5172
$ws.Range("I61").Value = 1597.5172
$ws.Range("K61").Value = 1597.5172
$ws.Range("M61").Value = -1385.5172
$ws.Range("H122").Value = 2622.7446
$ws.Range("I122").Value = 2035.871
$ws.Range("J122").Value = 3759.8125
$ws.Range("K122").Value = 6107.613
$ws.Range("L122").Value = 11279.4375
$ws.Range("M122").Value = -3657.613
$ws.Range("N122").Value = -16179.4375
$ws.Range("H136").Value = 1597.5172
$ws.Range("I136").Value = 1597.5172
$ws.Range("K136").Value = 4792.5516
$ws.Range("M136").Value = -2242.5516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1904.0526
$ws.Range("I105").Value = 1652.3214
$ws.Range("K105").Value = 1652.3214
$ws.Range("M105").Value = 94.67859999999996
$ws.Range("H134").Value = 5028.4443
$ws.Range("I134").Value = 4719.5
$ws.Range("K134").Value = 14158.5
$ws.Range("M134").Value = -11623.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3297.923
$ws.Range("I31").Value = 3082.8667
$ws.Range("K31").Value = 3082.8667
$ws.Range("M31").Value = -2787.8667
$ws.Range("H34").Value = 3297.923
$ws.Range("I34").Value = 3082.8667
$ws.Range("K34").Value = 3082.8667
$ws.Range("M34").Value = -2880.8667
$ws.Range("H104").Value = 47766
$ws.Range("I104").Value = 43299
$ws.Range("J104").Value = 49999.5
$ws.Range("K104").Value = 43299
$ws.Range("L104").Value = 49999.5
$ws.Range("M104").Value = -40678
$ws.Range("N104").Value = -55241.5
$ws.Range("H105").Value = 1072
$ws.Range("I105").Value = 1115.6
$ws.Range("K105").Value = 1115.6
$ws.Range("M105").Value = 631.4000000000001
$ws.Range("H122").Value = 2704.0625
$ws.Range("J122").Value = 2006.3334
$ws.Range("L122").Value = 6019.0002
$ws.Range("N122").Value = -10919.0002
$ws.Range("H132").Value = 529063.8
$ws.Range("I132").Value = 1648
$ws.Range("J132").Value = 1671798
$ws.Range("K132").Value = 4944
$ws.Range("L132").Value = 5015394
$ws.Range("M132").Value = -2414
$ws.Range("N132").Value = -5020454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 226.90909
$ws.Range("I2").Value = 111.5
$ws.Range("K2").Value = 669
$ws.Range("M2").Value = -556
$ws.Range("H6").Value = 107.94118
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 167.5
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 502.5
$ws.Range("M6").Value = -187
$ws.Range("N6").Value = -728.5
$ws.Range("H9").Value = 116909.5
$ws.Range("J9").Value = 2550
$ws.Range("L9").Value = 7650
$ws.Range("N9").Value = -8098
$ws.Range("H10").Value = 263.4
$ws.Range("I10").Value = 325.625
$ws.Range("J10").Value = 14.5
$ws.Range("K10").Value = 976.875
$ws.Range("L10").Value = 43.5
$ws.Range("M10").Value = -837.875
$ws.Range("N10").Value = -321.5
$ws.Range("H11").Value = 114101.77
$ws.Range("I11").Value = 30401.758
$ws.Range("J11").Value = 365201.8
$ws.Range("K11").Value = 91205.274
$ws.Range("L11").Value = 1095605.4
$ws.Range("M11").Value = -91065.274
$ws.Range("N11").Value = -1095885.4
$ws.Range("H117").Value = 1201.75
$ws.Range("J117").Value = 1495.1666
$ws.Range("L117").Value = 4485.4998
$ws.Range("N117").Value = -11369.4998
$ws.Range("H131").Value = 2969.3813
$ws.Range("J131").Value = 3127.3296
$ws.Range("L131").Value = 9381.988799999999
$ws.Range("N131").Value = -19461.9888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4989.231
$ws.Range("I132").Value = 4736.1
$ws.Range("K132").Value = 14208.3
$ws.Range("M132").Value = -11678.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1366980.1
$ws.Range("I40").Value = 2002411.4
$ws.Range("K40").Value = 2002411.4
$ws.Range("M40").Value = -2002275.4
$ws.Range("H55").Value = 649.6
$ws.Range("I55").Value = 217.45454
$ws.Range("J55").Value = 1838
$ws.Range("K55").Value = 217.45454
$ws.Range("L55").Value = 1838
$ws.Range("M55").Value = -44.45454000000001
$ws.Range("N55").Value = -2184
$ws.Range("H68").Value = 6391.0713
$ws.Range("I68").Value = 4696.5713
$ws.Range("J68").Value = 8085.5713
$ws.Range("K68").Value = 4696.5713
$ws.Range("L68").Value = 8085.5713
$ws.Range("M68").Value = -3947.5713
$ws.Range("N68").Value = -9583.5713
$ws.Range("H71").Value = 6391.0713
$ws.Range("I71").Value = 4696.5713
$ws.Range("J71").Value = 8085.5713
$ws.Range("K71").Value = 23482.8565
$ws.Range("L71").Value = 40427.85649999999
$ws.Range("M71").Value = -19738.8565
$ws.Range("N71").Value = -47915.85649999999
$ws.Range("H82").Value = 1399.2222
$ws.Range("I82").Value = 1172.8182
$ws.Range("K82").Value = 1172.8182
$ws.Range("M82").Value = -811.8181999999999
$ws.Range("H85").Value = 1399.2222
$ws.Range("I85").Value = 1172.8182
$ws.Range("K85").Value = 1172.8182
$ws.Range("M85").Value = 75.18180000000007
$ws.Range("H105").Value = 75299.5
$ws.Range("J105").Value = 75299.5
$ws.Range("L105").Value = 75299.5
$ws.Range("N105").Value = -82287.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3656418.5
$ws.Range("I62").Value = 8764105
$ws.Range("J62").Value = 8071.4287
$ws.Range("K62").Value = 8764105
$ws.Range("L62").Value = 8071.4287
$ws.Range("M62").Value = -8763481
$ws.Range("N62").Value = -9319.4287
$ws.Range("H65").Value = 3656418.5
$ws.Range("I65").Value = 8764105
$ws.Range("J65").Value = 8071.4287
$ws.Range("K65").Value = 43820525
$ws.Range("L65").Value = 40357.14350000001
$ws.Range("M65").Value = -43817405
$ws.Range("N65").Value = -46597.14350000001
$ws.Range("H126").Value = 1918.4286
$ws.Range("I126").Value = 1654.8334
$ws.Range("K126").Value = 4964.5002
$ws.Range("M126").Value = -2494.5002
$ws.Range("H132").Value = 449167.66
$ws.Range("I132").Value = 545064.9
$ws.Range("K132").Value = 1635194.7
$ws.Range("M132").Value = -1632664.7
$ws.Range("H136").Value = 5398.606
$ws.Range("I136").Value = 5764.2607
$ws.Range("K136").Value = 17292.7821
$ws.Range("M136").Value = -14742.7821
